# The edit permutes the per-row "price report" data (date, variety, quality,
# volume, prices, unit, origin, $/kg, kg/unit) across the 15 data rows
# (rows 2-16), while the descriptive/header-like columns (A,B,C,E,F,G,H,I,J)
# stay identical for every row (they already are identical across the sheet).
#
# Mapping: new row R gets the "D..T-subset" values that used to live in row
# $map[R] in the original sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$map = @{
    2  = 7
    3  = 12
    4  = 15
    5  = 16
    6  = 10
    7  = 9
    8  = 5
    9  = 2
    10 = 3
    11 = 4
    12 = 14
    13 = 8
    14 = 11
    15 = 6
    16 = 13
}

# Columns that actually change (snapshot + rewrite these only).
$cols = @(4, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20)  # D, K, L, M, N, O, P, Q, R, S, T

# 1) Snapshot all current values for rows 2..16 for the relevant columns,
#    so we can freely overwrite cells afterwards without losing data.
$snapshot = @{}
for ($r = 2; $r -le 16; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

# 2) Write back the values according to the permutation map.
for ($r = 2; $r -le 16; $r++) {
    $srcRow = $map[$r]
    $srcVals = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($r, $c).Value = $srcVals[$c]
    }
}
